# Sponsor logo/handle update: the "Happy Hour Sponsor" slide's Twitter
# handle placeholder reads "@aol" and needs to become "@aol_inc" (More
# updated sponsor logos).
$p = $ppt.ActivePresentation

$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $text = $shape.TextFrame.TextRange.Text
            if ($text -eq "@aol") {
                $targetSlide = $slide
                $targetShape = $shape
            }
        }
    }
}

if ($targetShape -eq $null) {
    # Fallback to the known location (slide 27, "Text Placeholder 2").
    $targetSlide = $p.Slides.Item(27)
    $targetShape = $targetSlide.Shapes.Item(3)
}

$tr = $targetShape.TextFrame.TextRange
$idx = $tr.Text.IndexOf("aol")
$run = $tr.Characters($idx + 1, 3)
$run.Text = "aol_inc"
